# Trading update: 2026-02-17 19:49:33
# Appends the newest trade (row 14) to both the "All Trades" and
# "MarketMaking" sheets. The trade is still open (no exit price / exit
# reason yet), so those columns are left blank.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 14

    $ws.Cells.Item($row, 1).Value  = 13
    # Leading apostrophe forces these to stay plain text instead of being
    # auto-parsed into date/time serial numbers.
    $ws.Cells.Item($row, 2).Value  = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value  = "'19:47:51"
    $ws.Cells.Item($row, 4).Value  = "MarketMaking"
    $ws.Cells.Item($row, 5).Value  = "DOWN"
    $ws.Cells.Item($row, 6).Value  = 0.48
    # Exit price: trade is still OPEN, so no exit price yet.
    $ws.Cells.Item($row, 8).Value  = "OPEN"
    $ws.Cells.Item($row, 9).Value  = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.2727272727273
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Exit reason: trade is still OPEN, so no exit reason yet.
    $ws.Cells.Item($row, 17).Value = 0
}
